# Update employment dates according to service records.
$d = $word.ActiveDocument

# 1) CSC / OBS job (bookmark "OBS"): "Юли 2003 – Октомври 2005" -> "Юли 2003 – Септември 2005"
$d.Content.Find.Execute("Юли 2003 – Октомври 2005", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Юли 2003 – Септември 2005", 2)

# 2) BIT job (bookmark "BIT"): "Септември 2001 – Март 2002" -> "Април 2001 – Февруари 2002"
$d.Content.Find.Execute("Септември 2001 – Март 2002", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Април 2001 – Февруари 2002", 2)

# 3) Total years of experience updated accordingly: "Повече от 13 години" -> "Повече от 14 години"
$d.Content.Find.Execute("Повече от 13 години", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Повече от 14 години", 2)

# 4) Last modified date text update: "7.06.2018 г." -> "26.11.2019 г."
$d.Content.Find.Execute("7.06.2018 г.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "26.11.2019 г.", 2)

# 5) The footer note ("Последна промяна: ..." and the following "latest version" line)
#    is now shown in italic, smaller (10pt) text.
$rng1 = $d.Content
$rng1.Find.Execute("Последна промяна:") | Out-Null
$rng1.Expand(4)

$rng2 = $d.Content
$rng2.Find.Execute("Най-пълна и текуща версия") | Out-Null
$rng2.Expand(4)

$noteRange = $d.Range($rng1.Start, $rng2.End)
$noteRange.Font.Italic = 1
$noteRange.Font.Size = 10
